$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Remove the old "Total" row (row 10) completely - it will be rebuilt lower
# down the sheet once the new order rows have been inserted.
# ---------------------------------------------------------------------------
$ws.Rows.Item(10).Delete() | Out-Null

# ---------------------------------------------------------------------------
# Update existing rows 2 & 3 - "Amazon Order" / "Hobby King Order" become
# "Amazon Order 1" / "Hobby King Order 1" (and their receipt file names).
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "Amazon Order 1"
$ws.Range("D2").Value = "Amazon Order 1.pdf"

$ws.Range("A3").Value = "Hobby King Order 1"
$ws.Range("D3").Value = "Hobby King Order 1.pdf"

# Rows 4-7 keep their existing data (Amazon Order 2/3/4, Hobby King Order 2).

# ---------------------------------------------------------------------------
# New order rows 8-11
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Amazon Order 5"
$ws.Range("B8").Value = "12/15/2017"
$ws.Range("C8").Value = "Brian"
$ws.Range("D8").Value = "Amazon Order 5.pdf"
$ws.Range("E8").Value = 149.87
$ws.Range("F8").Value = "Soldering Materials, Crimping tool, Silicone Wire, Alligator Clips, Heat Shrink Tubing, Resistors, Capacitors"

$ws.Range("A9").Value = "Digi-Key Order"
$ws.Range("B9").Value = "12/15/2017"
$ws.Range("C9").Value = "Brian"
$ws.Range("D9").Value = "Digi-Key Order.pdf"
$ws.Range("E9").Value = 45.31
$ws.Range("F9").Value = "Solder Wick, 0.1`" Connectors (Male and Female), Crimp Contacts"

$ws.Range("A10").Value = "Arrow Order"
$ws.Range("B10").Value = "12/15/2017"
$ws.Range("C10").Value = "Brian"
$ws.Range("D10").Value = "Arrow Order.pdf"
$ws.Range("E10").Value = 45.22
$ws.Range("F10").Value = "Slip Rings, Power MOSFETs, Rocker Switches"

$ws.Range("A11").Value = "Amazon Order 6*"
$ws.Range("B11").Value = "12/15/2017"
$ws.Range("C11").Value = "Brian"
$ws.Range("D11").Value = "Amazon Order 6.pdf"
$ws.Range("E11").Value = 147.68
$ws.Range("F11").Value = "Soldering Station, Tip Cleaner, Tip Tinner, Extra Iron Tips "

# Apply the same styles used by the existing data rows to the new rows.
$ws.Range("D8:D11").Style = "Hyperlink"
$ws.Range("E8:E11").NumberFormat = '"$"#,##0.00'

# ---------------------------------------------------------------------------
# Rebuild the "Total" row further down the sheet (row 13), summing all of
# the order rows above it, and add the funding note below (row 15).
# ---------------------------------------------------------------------------
$ws.Range("D13").Value = "Total"
$ws.Range("D13").Font.Bold = $true

$ws.Range("E13").Formula = "=SUM(E2:E12)"
$ws.Range("E13").NumberFormat = '"$"#,##0.00'

$ws.Range("A15").Value = "*Not sure if Feron will fund this order"

# ---------------------------------------------------------------------------
# Rebuild hyperlinks for every receipt cell in column D (the old collection
# is cleared first so stale links to renamed files are not left behind).
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("D2"), "Amazon Order 1.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "Hobby King Order 1.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D4"), "Amazon Order 2.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D5"), "Amazon Order 3.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D6"), "Amazon Order 4.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D7"), "Hobby King Order 2.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D8"), "Amazon Order 5.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D9"), "Digi-Key Order.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D10"), "Arrow Order.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D11"), "Amazon Order 6.pdf") | Out-Null

# Re-apply hyperlink styling so every receipt cell shares the same style
# index as before (Hyperlinks.Add can otherwise leave a stray style).
$ws.Range("D2:D11").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# Update the active selection to match the new end of the sheet.
# ---------------------------------------------------------------------------
$ws.Range("A16").Select()
